$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 4).Value = "35.324.69"
$ws.Cells.Item(2, 5).Value = "  +0.65%  "

$ws.Cells.Item(3, 4).Value = "1.881.30"
$ws.Cells.Item(3, 5).Value = "  -0.77%  "

$ws.Cells.Item(4, 5).Value = "  -0.65%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "245.99"
$ws.Cells.Item(5, 5).Value = "  -2.53%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.682"
$ws.Cells.Item(6, 5).Value = "  -2.24%  "

$ws.Cells.Item(7, 5).Value = "  -0.66%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "43.40"
$ws.Cells.Item(8, 5).Value = "  +4.61%  "

$ws.Cells.Item(9, 5).Value = "  +0.82%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "53.33"
$ws.Cells.Item(10, 5).Value = "  +2.08%  "

$ws.Cells.Item(11, 5).Value = "  -1.02%  "

$ws.Cells.Item(12, 5).Value = "  +0.06%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "13.54"
$ws.Cells.Item(13, 5).Value = "  +3.37%  "

$ws.Cells.Item(14, 4).Value = "2.153.74"
$ws.Cells.Item(14, 5).Value = "  -0.88%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.763"
$ws.Cells.Item(15, 5).Value = "  +5.08%  "

$ws.Cells.Item(17, 4).Value = "1.865.15"
$ws.Cells.Item(17, 5).Value = "  -2.22%  "

$ws.Cells.Item(18, 4).Value = "35.346.53"
$ws.Cells.Item(18, 5).Value = "  +0.66%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "73.64"
$ws.Cells.Item(19, 5).Value = "  -0.65%  "

$ws.Cells.Item(20, 5).Value = "  -1.14%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "244.50"
$ws.Cells.Item(21, 5).Value = "  -2.99%  "

$ws.Cells.Item(22, 5).Value = "  -0.95%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.20"
$ws.Cells.Item(23, 5).Value = "  +3.35%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.68"
$ws.Cells.Item(24, 5).Value = "  +10.45%  "

$ws.Cells.Item(25, 5).Value = "  -0.62%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.17"
$ws.Cells.Item(26, 5).Value = "  -4.56%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "164.85"
$ws.Cells.Item(27, 5).Value = "  -2.11%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "8.71"
$ws.Cells.Item(28, 5).Value = "  +1.86%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "18.30"
$ws.Cells.Item(29, 5).Value = "  -0.60%  "

$ws.Cells.Item(30, 5).Value = "  -0.92%  "

$ws.Cells.Item(31, 5).Value = "  -0.30%  "

$ws.Cells.Item(32, 5).Value = "  +0.27%  "

$ws.Cells.Item(33, 5).Value = "  -1.09%  "

$ws.Cells.Item(34, 5).Value = "  -0.66%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.82"
$ws.Cells.Item(35, 5).Value = "  -3.12%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.44"
$ws.Cells.Item(36, 5).Value = "  -11.30%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.852"
$ws.Cells.Item(37, 5).Value = "  +0.51%  "

$ws.Cells.Item(38, 5).Value = "  -2.68%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.0726"
$ws.Cells.Item(39, 5).Value = "  +10.14%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "17.40"
$ws.Cells.Item(40, 5).Value = "  -0.69%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0219"
$ws.Cells.Item(41, 5).Value = "  +2.71%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "97.12"
$ws.Cells.Item(42, 5).Value = "  -1.63%  "

$ws.Cells.Item(43, 5).Value = "  -2.21%  "

$ws.Cells.Item(44, 5).Value = "  +1.31%  "

$ws.Cells.Item(45, 4).Value = "1.309.44"
$ws.Cells.Item(45, 5).Value = "  +0.59%  "

$ws.Cells.Item(46, 5).Value = "  +4.04%  "

$ws.Cells.Item(47, 5).Value = "  -1.09%  "

$ws.Cells.Item(48, 5).Value = "  -0.25%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "11.96"
$ws.Cells.Item(49, 5).Value = "  +0.28%  "

$ws.Cells.Item(50, 5).Value = "  -3.93%  "

$ws.Cells.Item(51, 5).Value = "  -1.45%  "
